$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated "K" (strikeouts) column values replacing the old "Strike#" derived
# values. Column G corresponds to header "K" (row 1).
$kValues = @{
    2  = 2
    3  = 1
    4  = 1
    5  = 0
    6  = 1
    7  = 1
    8  = 0
    9  = 2
    10 = 2
    11 = 1
    12 = 0
    13 = 2
    14 = 2
    15 = 0
    16 = 2
    17 = 1
    18 = 2
    19 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
